$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '27.521.81'
$c.ClearFormats()
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  +4.24%  '
$c.ClearFormats()
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.841.78'
$c.ClearFormats()
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  +3.73%  '
$c.ClearFormats()
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.029'
$c.ClearFormats()
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  +2.74%  '
$c.ClearFormats()
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '318.60'
$c.ClearFormats()
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  +4.03%  '
$c.ClearFormats()
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '1.026'
$c.ClearFormats()
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  +2.38%  '
$c.ClearFormats()
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.4373'
$c.ClearFormats()
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  +3.33%  '
$c.ClearFormats()
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.3731'
$c.ClearFormats()
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  +3.30%  '
$c.ClearFormats()
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.07387'
$c.ClearFormats()
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +3.31%  '
$c.ClearFormats()
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.8752'
$c.ClearFormats()
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  +4.37%  '
$c.ClearFormats()
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '21.46'
$c.ClearFormats()
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +5.06%  '
$c.ClearFormats()
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '1.851.23'
$c.ClearFormats()
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  +4.42%  '
$c.ClearFormats()
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '5.494'
$c.ClearFormats()
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +4.53%  '
$c.ClearFormats()
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '6.685'
$c.ClearFormats()
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  +3.62%  '
$c.ClearFormats()
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.07137'
$c.ClearFormats()
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  +3.39%  '
$c.ClearFormats()
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '82.69'
$c.ClearFormats()
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  +4.65%  '
$c.ClearFormats()
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '1.030'
$c.ClearFormats()
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  +2.83%  '
$c.ClearFormats()
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '0.000009002'
$c.ClearFormats()
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  +3.17%  '
$c.ClearFormats()
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '1.025'
$c.ClearFormats()
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  +2.34%  '
$c.ClearFormats()
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '15.40'
$c.ClearFormats()
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  +3.28%  '
$c.ClearFormats()
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '27.530.45'
$c.ClearFormats()
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  +4.24%  '
$c.ClearFormats()
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '5.230'
$c.ClearFormats()
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  +2.57%  '
$c.ClearFormats()
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '11.19'
$c.ClearFormats()
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  +2.51%  '
$c.ClearFormats()
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.067.18'
$c.ClearFormats()
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  +3.95%  '
$c.ClearFormats()
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '156.85'
$c.ClearFormats()
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +3.45%  '
$c.ClearFormats()
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '1.920'
$c.ClearFormats()
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  +7.05%  '
$c.ClearFormats()
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '18.70'
$c.ClearFormats()
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  +3.61%  '
$c.ClearFormats()
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '5.248'
$c.ClearFormats()
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  +3.54%  '
$c.ClearFormats()
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.935'
$c.ClearFormats()
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  +5.09%  '
$c.ClearFormats()
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '116.28'
$c.ClearFormats()
$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  +1.81%  '
$c.ClearFormats()
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.09073'
$c.ClearFormats()
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +2.80%  '
$c.ClearFormats()
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  +7.86%  '
$c.ClearFormats()
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.7645'
$c.ClearFormats()
$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  +4.98%  '
$c.ClearFormats()
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '4.489'
$c.ClearFormats()
$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  +3.94%  '
$c.ClearFormats()
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  +5.13%  '
$c.ClearFormats()
$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  +2.66%  '
$c.ClearFormats()
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '1.147'
$c.ClearFormats()
$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  +5.01%  '
$c.ClearFormats()
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.01969'
$c.ClearFormats()
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  +4.37%  '
$c.ClearFormats()
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.05253'
$c.ClearFormats()
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  +2.74%  '
$c.ClearFormats()
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.5171'
$c.ClearFormats()
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +4.97%  '
$c.ClearFormats()
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.793'
$c.ClearFormats()
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  +7.41%  '
$c.ClearFormats()
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.1665'
$c.ClearFormats()
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  +3.50%  '
$c.ClearFormats()
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '6.619'
$c.ClearFormats()
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  +4.54%  '
$c.ClearFormats()
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '8.537'
$c.ClearFormats()
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  +5.54%  '
$c.ClearFormats()
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '109.27'
$c.ClearFormats()
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  +4.35%  '
$c.ClearFormats()
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '10.61'
$c.ClearFormats()
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  +3.50%  '
$c.ClearFormats()
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  +2.63%  '
$c.ClearFormats()
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.706'
$c.ClearFormats()
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  +4.99%  '
$c.ClearFormats()
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.4642'
$c.ClearFormats()
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  +4.43%  '
$c.ClearFormats()
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '1.895'
$c.ClearFormats()
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  +10.64%  '
$c.ClearFormats()
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.06330'
$c.ClearFormats()
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  +2.55%  '
$c.ClearFormats()
